# BUGFIX: Escape characters in IQ_bullet_list
# Updates the "script_IQSlidedeck.R | Date: ..." timestamp text that
# appears in the footer placeholder of several slides.

$p = $ppt.ActivePresentation

$updates = @{
    3  = "script_IQSlidedeck.R | Date: 2021-03-09 17:50:17"
    5  = "script_IQSlidedeck.R | Date: 2021-03-09 17:50:17"
    6  = "script_IQSlidedeck.R | Date: 2021-03-09 17:50:21"
    7  = "script_IQSlidedeck.R | Date: 2021-03-09 17:50:22"
    8  = "script_IQSlidedeck.R | Date: 2021-03-09 17:50:28"
    10 = "script_IQSlidedeck.R | Date: 2021-03-09 17:50:18"
    11 = "script_IQSlidedeck.R | Date: 2021-03-09 17:50:19"
    12 = "script_IQSlidedeck.R | Date: 2021-03-09 17:50:20"
}

foreach ($slideIdx in $updates.Keys) {
    $s = $p.Slides.Item($slideIdx)
    $footer = $s.Shapes.Item($s.Shapes.Count)
    $footer.TextFrame.TextRange.Text = $updates[$slideIdx]
}
